$wb = $excel.ActiveWorkbook

# --- Sheet "Create WO" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Create WO")

# New header cells H1:L1
$ws1.Range("H1").Value = "Background Processing"
$ws1.Range("I1").Value = "Labor Booking User"
$ws1.Range("J1").Value = "SiteID"
$ws1.Range("K1").Value = "Location ID"
$ws1.Range("L1").Value = "Location Number"
$ws1.Range("K1").Style = "Normal"

# Row 2 new values
$ws1.Range("H2").Value = $true
$ws1.Range("I2").Value = "a9L5f000000se6a"
$ws1.Range("J2").Value = "a9A5f000000YM1d"
$ws1.Range("K2").Value = "a8s5f000000PtAf"
$ws1.Range("L2").Value = "SY_ReceiptLoc"
$ws1.Range("K2").Style = "Normal"

# Row 3 new values, plus removal of old E3/F3
$ws1.Range("E3").ClearContents()
$ws1.Range("F3").ClearContents()
$ws1.Range("H3").Value = $false
$ws1.Range("I3").Value = "a9L5f000000se6a"
$ws1.Range("J3").Value = "a9A5f000000YM1d"
$ws1.Range("K3").Value = "a8s5f000000PtAf"
$ws1.Range("L3").Value = "SY_ReceiptLoc"
$ws1.Range("K3").Style = "Normal"

# --- Sheet "Add Component" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Add Component")
$ws2.Range("G4").ClearContents()
$ws2.Range("H4").ClearContents()
$ws2.Range("G7").ClearContents()
$ws2.Range("H7").ClearContents()
$ws2.Range("D8").Select()

# --- Selections / active sheet ---
$ws1.Activate()
$ws1.Range("F8").Select()
